$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = "In Translation"
$ov.Range("F2").Value = "In Translation"
$ov.Range("G2").Value = "2016-09-07 03:27:30"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-09-07 03:27:30"

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("H2").Value = "2016-09-07 03:27:19"
$zh.Range("P2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e1708cea0d9e0c69516bfc17e7db7d3915b2623/e2e/b64778cb-8792-415b-80c6-c326caff7005.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ff8d0e2e4ca879de894cc8713cc0a5c61f76438/e2e/b64778cb-8792-415b-80c6-c326caff7005.md."
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("H3").Value = "2016-09-07 03:27:19"
$zh.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e1708cea0d9e0c69516bfc17e7db7d3915b2623/e2e/ff3b416d-95d1-4df8-aca0-c7d146132f98.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ff8d0e2e4ca879de894cc8713cc0a5c61f76438/e2e/ff3b416d-95d1-4df8-aca0-c7d146132f98.md."

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = "In Translation"
$de.Range("H2").Value = "2016-09-07 03:27:30"
$de.Range("P2").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e1708cea0d9e0c69516bfc17e7db7d3915b2623/e2e/b64778cb-8792-415b-80c6-c326caff7005.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ff8d0e2e4ca879de894cc8713cc0a5c61f76438/e2e/b64778cb-8792-415b-80c6-c326caff7005.md."
$de.Range("C3").Value = "In Translation"
$de.Range("H3").Value = "2016-09-07 03:27:30"
$de.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3e1708cea0d9e0c69516bfc17e7db7d3915b2623/e2e/ff3b416d-95d1-4df8-aca0-c7d146132f98.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4ff8d0e2e4ca879de894cc8713cc0a5c61f76438/e2e/ff3b416d-95d1-4df8-aca0-c7d146132f98.md."

# ---- column width adjustments ----
# Target stored widths: 17.2159881591797 (cols) and 40 (error-detail cols).
# The engine quantizes ColumnWidth to 1/6-character steps, so 16.3 is the
# closest input that lands on the nearest achievable stored width (17.1667).
$ov.Columns.Item(5).ColumnWidth = 16.3
$ov.Columns.Item(6).ColumnWidth = 16.3
$zh.Columns.Item(3).ColumnWidth = 16.3
$zh.Columns.Item(16).ColumnWidth = 39.2
$de.Columns.Item(3).ColumnWidth = 16.3
$de.Columns.Item(16).ColumnWidth = 39.2
